# Update Excalibur_Profits market-price figures (scheduled data refresh).
# Each cell below is a raw (non-formula) numeric value pulled from the
# latest market snapshot; only the H/I/J/K/L/M/N "price/profit" columns move.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H11").Value = 106.4
$ws.Range("I11").Value = 106.4
$ws.Range("K11").Value = 106.4
$ws.Range("M11").Value = 33.59999999999999
$ws.Range("H32").Value = 1381.5714
$ws.Range("J32").Value = 1356.3334
$ws.Range("L32").Value = 1356.3334
$ws.Range("N32").Value = -2008.3334
$ws.Range("H33").Value = 204.875
$ws.Range("I33").Value = 204.8125
$ws.Range("K33").Value = 204.8125
$ws.Range("M33").Value = 24.1875
$ws.Range("H40").Value = 2593.8333
$ws.Range("I40").Value = 2323.25
$ws.Range("K40").Value = 2323.25
$ws.Range("M40").Value = -2148.25
$ws.Range("H69").Value = 8802.842000000001
$ws.Range("I69").Value = 5399.6665
$ws.Range("K69").Value = 16198.9995
$ws.Range("M69").Value = -15324.9995
$ws.Range("H72").Value = 8802.842000000001
$ws.Range("I72").Value = 5399.6665
$ws.Range("K72").Value = 48596.9985
$ws.Range("M72").Value = -44228.9985
$ws.Range("H74").Value = 7510.2104
$ws.Range("I74").Value = 4929
$ws.Range("J74").Value = 9015.916999999999
$ws.Range("K74").Value = 4929
$ws.Range("L74").Value = 9015.916999999999
$ws.Range("M74").Value = -3993
$ws.Range("N74").Value = -10887.917
$ws.Range("H77").Value = 7510.2104
$ws.Range("I77").Value = 4929
$ws.Range("J77").Value = 9015.916999999999
$ws.Range("K77").Value = 24645
$ws.Range("L77").Value = 45079.585
$ws.Range("M77").Value = -19965
$ws.Range("N77").Value = -54439.585
$ws.Range("H101").Value = 2159.7334
$ws.Range("J101").Value = 4018.6
$ws.Range("L101").Value = 12055.8
$ws.Range("N101").Value = -15299.8
$ws.Range("H132").Value = 1735.3611
$ws.Range("I132").Value = 1735.3611
$ws.Range("K132").Value = 5206.0833
$ws.Range("M132").Value = -2676.0833
$ws.Range("H138").Value = 2730.6482
$ws.Range("I138").Value = 1906.92
$ws.Range("J138").Value = 3440.7585
$ws.Range("K138").Value = 5720.76
$ws.Range("L138").Value = 10322.2755
$ws.Range("M138").Value = -580.7600000000002
$ws.Range("N138").Value = -20602.2755

# --- Sheet: ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H88").Value = 1482.0526
$ws.Range("I88").Value = 1395.8182
$ws.Range("J88").Value = 1600.625
$ws.Range("K88").Value = 1395.8182
$ws.Range("L88").Value = 1600.625
$ws.Range("M88").Value = -989.8181999999999
$ws.Range("N88").Value = -2412.625
$ws.Range("H91").Value = 1482.0526
$ws.Range("I91").Value = 1395.8182
$ws.Range("J91").Value = 1600.625
$ws.Range("K91").Value = 1395.8182
$ws.Range("L91").Value = 1600.625
$ws.Range("M91").Value = 8.181800000000067
$ws.Range("N91").Value = -4408.625
$ws.Range("H97").Value = 1740.7778
$ws.Range("I97").Value = 1715.6666
$ws.Range("J97").Value = 1866.3334
$ws.Range("K97").Value = 1715.6666
$ws.Range("L97").Value = 1866.3334
$ws.Range("M97").Value = -1219.6666
$ws.Range("N97").Value = -2858.3334
$ws.Range("H110").Value = 1136.6
$ws.Range("I110").Value = 1040.6666
$ws.Range("K110").Value = 1040.6666
$ws.Range("M110").Value = 1004.3334

# --- Sheet: BSM (index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 2544
$ws.Range("J20").Value = 2225.8333
$ws.Range("L20").Value = 2225.8333
$ws.Range("N20").Value = -2719.8333
$ws.Range("H86").Value = 2409.5
$ws.Range("I86").Value = 2900.8
$ws.Range("K86").Value = 2900.8
$ws.Range("M86").Value = -1777.8
$ws.Range("H89").Value = 2409.5
$ws.Range("I89").Value = 2900.8
$ws.Range("K89").Value = 14504
$ws.Range("M89").Value = -8888
$ws.Range("H94").Value = 1172.4584
$ws.Range("I94").Value = 1129.1052
$ws.Range("J94").Value = 1337.2
$ws.Range("K94").Value = 1129.1052
$ws.Range("L94").Value = 1337.2
$ws.Range("M94").Value = -678.1052
$ws.Range("N94").Value = -2239.2
$ws.Range("H105").Value = 2565.75
$ws.Range("I105").Value = 2378.9
$ws.Range("K105").Value = 2378.9
$ws.Range("M105").Value = -631.9000000000001

# --- Sheet: CRP (index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H17").Value = 18750
$ws.Range("I17").Value = 12500
$ws.Range("J17").Value = 25000
$ws.Range("K17").Value = 12500
$ws.Range("L17").Value = 25000
$ws.Range("M17").Value = -12326
$ws.Range("N17").Value = -25348
$ws.Range("H58").Value = 1793.3112
$ws.Range("I58").Value = 1391.7646
$ws.Range("J58").Value = 3034.4546
$ws.Range("K58").Value = 1391.7646
$ws.Range("L58").Value = 3034.4546
$ws.Range("M58").Value = -1188.7646
$ws.Range("N58").Value = -3440.4546
$ws.Range("H132").Value = 6671.864
$ws.Range("I132").Value = 7246.1665
$ws.Range("K132").Value = 21738.4995
$ws.Range("M132").Value = -19208.4995
$ws.Range("H134").Value = 1321.7333
$ws.Range("I134").Value = 1198.1923
$ws.Range("K134").Value = 3594.5769
$ws.Range("M134").Value = -1059.5769
$ws.Range("H136").Value = 1793.3112
$ws.Range("I136").Value = 1391.7646
$ws.Range("J136").Value = 3034.4546
$ws.Range("K136").Value = 4175.293799999999
$ws.Range("L136").Value = 9103.363799999999
$ws.Range("M136").Value = -1625.293799999999
$ws.Range("N136").Value = -14203.3638
$ws.Range("H141").Value = 253384.92
$ws.Range("J141").Value = 253384.92
$ws.Range("L141").Value = 253384.92
$ws.Range("N141").Value = -263744.92

# --- Sheet: CUL (index 5) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H14").Value = 550.94116
$ws.Range("I14").Value = 550.94116
$ws.Range("K14").Value = 1652.82348
$ws.Range("M14").Value = -1479.82348
$ws.Range("H81").Value = 5264.857
$ws.Range("I81").Value = 899.5
$ws.Range("J81").Value = 7011
$ws.Range("K81").Value = 2698.5
$ws.Range("L81").Value = 21033
$ws.Range("M81").Value = -1575.5
$ws.Range("N81").Value = -23279
$ws.Range("H84").Value = 5264.857
$ws.Range("I84").Value = 899.5
$ws.Range("J84").Value = 7011
$ws.Range("K84").Value = 8095.5
$ws.Range("L84").Value = 63099
$ws.Range("M84").Value = -2479.5
$ws.Range("N84").Value = -74331
$ws.Range("H92").Value = 784.36365
$ws.Range("I92").Value = 535.5
$ws.Range("K92").Value = 1606.5
$ws.Range("M92").Value = -358.5
$ws.Range("H131").Value = 1474.069
$ws.Range("J131").Value = 1649.9565
$ws.Range("L131").Value = 4949.8695
$ws.Range("N131").Value = -15029.8695

# --- Sheet: GSM (index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H132").Value = 10503.131
$ws.Range("I132").Value = 10172.454
$ws.Range("K132").Value = 30517.362
$ws.Range("M132").Value = -27987.362

# --- Sheet: LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 2255.5264
$ws.Range("I7").Value = 2153
$ws.Range("K7").Value = 2153
$ws.Range("M7").Value = -2041
$ws.Range("H22").Value = 115017.664
$ws.Range("I22").Value = 505000.5
$ws.Range("J22").Value = 3594
$ws.Range("K22").Value = 505000.5
$ws.Range("L22").Value = 3594
$ws.Range("M22").Value = -504705.5
$ws.Range("N22").Value = -4184
$ws.Range("H27").Value = 115017.664
$ws.Range("I27").Value = 505000.5
$ws.Range("J27").Value = 3594
$ws.Range("K27").Value = 505000.5
$ws.Range("L27").Value = 3594
$ws.Range("M27").Value = -504893.5
$ws.Range("N27").Value = -3808
$ws.Range("H40").Value = 5636
$ws.Range("I40").Value = 4896.533
$ws.Range("K40").Value = 4896.533
$ws.Range("M40").Value = -4760.533
$ws.Range("H46").Value = 1515.8636
$ws.Range("I46").Value = 1186.4
$ws.Range("J46").Value = 2221.8572
$ws.Range("K46").Value = 1186.4
$ws.Range("L46").Value = 2221.8572
$ws.Range("M46").Value = -998.4000000000001
$ws.Range("N46").Value = -2597.8572
$ws.Range("H55").Value = 1724.6428
$ws.Range("I55").Value = 439.66666
$ws.Range("J55").Value = 2688.375
$ws.Range("K55").Value = 439.66666
$ws.Range("L55").Value = 2688.375
$ws.Range("M55").Value = -266.66666
$ws.Range("N55").Value = -3034.375
$ws.Range("H82").Value = 2224.889
$ws.Range("I82").Value = 2176
$ws.Range("J82").Value = 2249.3333
$ws.Range("K82").Value = 2176
$ws.Range("L82").Value = 2249.3333
$ws.Range("M82").Value = -1815
$ws.Range("N82").Value = -2971.3333
$ws.Range("H85").Value = 2224.889
$ws.Range("I85").Value = 2176
$ws.Range("J85").Value = 2249.3333
$ws.Range("K85").Value = 2176
$ws.Range("L85").Value = 2249.3333
$ws.Range("M85").Value = -928
$ws.Range("N85").Value = -4745.3333
$ws.Range("H107").Value = 19499.5
$ws.Range("I107").Value = 19499.5
$ws.Range("K107").Value = 19499.5
$ws.Range("M107").Value = -17579.5
$ws.Range("H126").Value = 2255.5264
$ws.Range("I126").Value = 2153
$ws.Range("K126").Value = 6459
$ws.Range("M126").Value = -3989
$ws.Range("H136").Value = 6590.636
$ws.Range("I136").Value = 5834.5654
$ws.Range("K136").Value = 17503.6962
$ws.Range("M136").Value = -14953.6962

# --- Sheet: WVR (index 8) ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H126").Value = 4216.6665
$ws.Range("I126").Value = 4216.6665
$ws.Range("K126").Value = 12649.9995
$ws.Range("M126").Value = -10179.9995
$ws.Range("H132").Value = 3571.8333
$ws.Range("I132").Value = 2797.7354
$ws.Range("K132").Value = 8393.206200000001
$ws.Range("M132").Value = -5863.206200000001
